# 1. UI_btn 추가 중 멈춤.
# Update row 13 (플레이어 스킬 추가(다른 종류)) with remark + start date/time + end time.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# E13: remark text (plain style, matches default style already used in column)
$ws.Range("E13").Value = "스킬 2개 추가. (즉발형)"

# F13: start date -> copy number format from an existing date cell (F4) so the
# same style (s="2") is reused instead of creating a brand new style entry.
$ws.Range("F4").Copy()
$ws.Range("F13").PasteSpecial(-4122)
$ws.Range("F13").Value = 44838

# G13: start time -> copy number format from an existing time cell (G4) so the
# same style (s="3") is reused.
$ws.Range("G4").Copy()
$ws.Range("G13").PasteSpecial(-4122)
$ws.Range("G13").Value = 0.68402777777777779

# I13: end time -> same time style as above.
$ws.Range("I4").Copy()
$ws.Range("I13").PasteSpecial(-4122)
$ws.Range("I13").Value = 0.76527777777777783

# New rows 17 and 18: additional "사운드" tasks.
$ws.Range("B17").Value = 10
$ws.Range("C17").Value = "사운드"
$ws.Range("D17").Value = "맵 사운드 출력"

$ws.Range("C18").Value = "사운드"
$ws.Range("D18").Value = "타격 사운드 출력"

# Update selection to match the author's final cursor position.
$ws.Range("E18").Select()
